# Add "Merge Sorted Array" (LeetCode 88) as a new row (row 22) on the
# LeetCode problem tracker sheet, matching the existing table's layout:
#   A: Name            B: Category   C: Blind75 problem?
#   D: Revisit         E: Difficulty F: Relative Difficulty
#   G: Local Path (hyperlink to the local solution file)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "Merge Sorted Array"
$ws.Range("B22").Value = "Array"
$ws.Range("C22").Value = "No"
$ws.Range("D22").Value = "Yes"
$ws.Range("E22").Value = "Easy"
$ws.Range("F22").Value = "Medium"
$ws.Range("G22").Value = "88 - Merge Sorted Array"

# Hyperlink G22 to the local solution file, same pattern as the rows above it.
$ws.Hyperlinks.Add($ws.Range("G22"), "88 - Merge Sorted Array")

# Keep G22 styled like the other "Local Path" hyperlink cells in the column.
$ws.Range("G22").Style = $ws.Range("G21").Style

# Match the selection left behind in the saved workbook.
$ws.Range("M23").Select()
